$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # "总计" summary sheet
$ws3 = $wb.Worksheets.Item(2)          # currently "2022-Q2"; will become "2022-Q3"

# ---------------------------------------------------------------------------
# 1) Duplicate the existing "2022-Q2" sheet so its original fund-holding data
#    survives unmodified under a brand-new sheet placed right after it
#    (this becomes sheetId=3 / the new sheet3.xml).
# ---------------------------------------------------------------------------
$ws3.Copy($null, $ws3)

# 2) Rename: the original sheet becomes "2022-Q3"; the fresh copy keeps the
#    "2022-Q2" name (so the Q2 data keeps living under the "2022-Q2" tab).
$ws3.Name = "2022-Q3"
$wsQ2 = $wb.Worksheets.Item(3)
$wsQ2.Name = "2022-Q2"

# ---------------------------------------------------------------------------
# 3) Replace the renamed sheet's (old Q2) fund-holding data with the new
#    2022-Q3 figures.
# ---------------------------------------------------------------------------
$ws3.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    # Pull formatting from the "总计" sheet's header cell (style used by newer data)
    # then overwrite with the real header text.
    $ws1.Range("B1").Copy($ws3.Cells.Item(1, $col))
    $ws3.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$rows = @(
    @(0, "005774", "华夏产业升级混合A", "24.29", "93.85", "5.59", "1.3578", 10),
    @(1, "501079", "大成科创主题混合（LOF）A", "10.00", "85.01", "8.66", "0.8660", 4),
    @(2, "015059", "华夏产业升级混合C", "8.92", "93.85", "5.59", "0.4986", 10),
    @(3, "012473", "大成成长回报六个月持有混合A", "7.48", "80.86", "5.66", "0.4234", 6),
    @(4, "160926", "大成创业板两年定期开放混合A", "7.68", "64.14", "3.92", "0.3011", 5),
    @(5, "010371", "大成成长进取混合A", "3.61", "80.75", "5.40", "0.1949", 5),
    @(6, "009798", "大成创业板两年定期开放混合C", "2.71", "64.14", "3.92", "0.1062", 5),
    @(7, "010372", "大成成长进取混合C", "1.52", "80.75", "5.40", "0.0821", 5),
    @(8, "012474", "大成成长回报六个月持有混合C", "0.37", "80.86", "5.66", "0.0209", 6),
    @(9, "016198", "大成科创主题混合（LOF）C", "0.01", "85.01", "8.66", "0.0009", 4)
)

foreach ($row in $rows) {
    $r = [int]$row[0] + 2

    # Column A: row index, styled like the "总计" sheet's index column.
    $ws1.Range("A2").Copy($ws3.Cells.Item($r, 1))
    $ws3.Cells.Item($r, 1).Value = $row[0]

    # Columns B-G: free-text values (fund code / name / percentages). Force
    # text storage so strings that look numeric (e.g. the leading-zero fund
    # code "005774") aren't silently coerced into numbers.
    $textRange = $ws3.Range($ws3.Cells.Item($r, 2), $ws3.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"
    for ($c = 2; $c -le 7; $c++) {
        $ws3.Cells.Item($r, $c).Value = [string]$row[$c - 1]
    }
    $textRange.Style = "Normal"

    # Column H: numeric position-ranking.
    $ws3.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------------
# 4) Update the "总计" summary sheet: insert the new 2022-Q3 totals as row 2
#    and push the existing 2022-Q2 totals down to row 3.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(2).Copy($ws1.Rows.Item(3))

$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 10
$ws1.Range("D2").Value = 3.85
